$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Final_Matches"

# Header row (bold, centered, bordered to match the workbook's existing header style)
$ws.Range("A1").Value = "AZ.CT/LABEL"
$ws.Range("B1").Value = "ASCTB.CT/LABEL"
$hdr = $ws.Range("A1:B1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

$data = @(
    @("monocyte", "monocyte"),
    @("memory B cell", "memory B cell"),
    @("naive B cell", "naive B cell"),
    @("regulatory T cell", "regulatory T cell"),
    @("naive thymus-derived CD4-positive, alpha-beta T cell", "naive thymus-derived CD4-positive, alpha-beta T cell"),
    @("naive thymus-derived CD8-positive, alpha-beta T cell", "naive thymus-derived CD8-positive, alpha-beta T cell"),
    @("CD16-negative, CD56-bright natural killer cell, human", "CD16-negative, CD56-bright natural killer cell, human"),
    @("effector CD8-positive, alpha-beta T cell", "effector CD8-positive, alpha-beta T cell"),
    @("CD141-positive myeloid dendritic cell", "CD141-positive myeloid dendritic cell"),
    @("CD14-low, CD16-positive monocyte", "CD14-low, CD16-positive monocyte"),
    @("CD1c-positive myeloid dendritic cell", "CD1c-positive myeloid dendritic cell"),
    @("mature B cell", "lymphocyte of B lineage"),
    @("CD4-positive, alpha-beta T cell", "T cell"),
    @("CD8-positive, alpha-beta T cell", "T cell"),
    @("mature T cell", "T cell"),
    @("CD4-positive, alpha-beta cytotoxic T cell", "T cell"),
    @("activated CD4-positive, alpha-beta T cell", "T cell"),
    @("CD4-positive, alpha-beta memory T cell", "T cell"),
    @("effector CD4-positive, alpha-beta T cell", "T cell"),
    @("activated CD8-positive, alpha-beta T cell", "T cell"),
    @("CD8-positive, alpha-beta memory T cell", "T cell"),
    @("CD4-negative, CD8-negative, alpha-beta intraepithelial T cell", "T cell"),
    @("gamma-delta T cell", "T cell"),
    @("mucosal invariant T cell", "T cell"),
    @("CD16-positive, CD56-dim natural killer cell, human", "mature natural killer cell"),
    @("plasmablast", "lymphocyte of B lineage"),
    @("plasma cell", "lymphocyte of B lineage"),
    @("plasmacytoid dendritic cell", "plasmacytoid dendritic cell, human"),
    @("natural killer cell", "mature natural killer cell")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
